# Agrego de nuevo el Excel de localidades
#
# The "localidad" column (A) previously stored several abbreviated place
# names ("Pto. Bandera", "Gdor. Gregores", "Pto. Santa Cruz",
# "Cmte.L.Piedrabuena", "Pto. San Julian"). They are expanded back to their
# full names here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value  = "Puerto Bandera"
$ws.Range("A12").Value = "Gobernador Gregores"
$ws.Range("A18").Value = "Puerto Santa Cruz"
$ws.Range("A19").Value = "Comandante Luis Piedrabuena"
$ws.Range("A20").Value = "Puerto San Julian"

# Column A needs to be a bit wider to fit the longer names.
$ws.Columns.Item(1).ColumnWidth = 25

# Leave the selection where the author left it when saving.
$ws.Range("D24").Select() | Out-Null
